$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the shared header text "Kraken2" -> "Kraken 2"
$ws.Range("Z2").Value = "Kraken 2"

# Replace the static values in Z3:Z25 with a formula referencing the header cell
for ($r = 3; $r -le 25; $r++) {
    $ws.Range("Z$r").Formula = "=Z`$2"
}

# Update the active selection to match the new edit range
$ws.Range("Z3:Z25").Select()
